# "Add ch9 report & modified ch6"
#
# performance.xlsx (Sheet1) records a speedup / Karp-Flatt benchmark.
# The "Real execution time" measurements (row 2, B2:I2) were re-run and
# the new timings are written in; every dependent formula (row 3 estimate,
# row 4 speedup, row 5 Karp-Flatt metric, row 7 lambda) recalculates off
# of them automatically. The sheet's active selection and the workbook
# window state were also left in a different UI state by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-measured "Real execution time" values (row 2, columns B:I) ---
$newValues = @(
    2.1521389480000002,
    1.085650921,
    0.72222399699999995,
    0.53893685300000005,
    0.434614897,
    0.372735977,
    0.31239581100000002,
    0.27757406200000001
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $col = 2 + $i   # column B = 2 ... I = 9
    $ws.Cells.Item(2, $col).Value = $newValues[$i]
}

# Recalculate the whole workbook so the dependent formulas (estimate time,
# speedup, Karp-Flatt metric, lambda) pick up the new measurements.
$excel.CalculateFullRebuild()

# --- Leave the sheet's selection on A1:I5 (whole data block) ---
$ws.Activate()
$ws.Range("A1:I5").Select()

# --- Workbook window was left minimized ---
$excel.ActiveWindow.WindowState = -4140   # xlMinimized
